# Append the four new comment rows (52-55) that were missing, fixing the
# invalid date shown for comment "fdgrt" in Firefox (04-14-2023 needs to be
# stored as text, not auto-converted to a date serial number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54 -> id 52, "fdgrt", 04-14-2023
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = "fdgrt"
# Force text storage for the date-like string so Excel doesn't silently
# convert it into a date serial number; then drop back to the default
# ("Normal") style so the cell keeps the same look as its neighbours.
$ws.Range("C54").NumberFormat = "@"
$ws.Range("C54").Value = "04-14-2023"
$ws.Range("C54").Style = "Normal"

# Row 55 -> id 53, "sawesda", ISO timestamp
$ws.Range("A55").Value = 53
$ws.Range("B55").Value = "sawesda"
$ws.Range("C55").Value = "2023-04-15T01:50:48.080Z"

# Row 56 -> id 54, "ersdfs", ISO timestamp
$ws.Range("A56").Value = 54
$ws.Range("B56").Value = "ersdfs"
$ws.Range("C56").Value = "2023-04-15T01:52:14.613Z"

# Row 57 -> id 55, "qwewqe", ISO timestamp
$ws.Range("A57").Value = 55
$ws.Range("B57").Value = "qwewqe"
$ws.Range("C57").Value = "2023-04-15T01:52:44.738Z"
